$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 = copy row 44 formatting exactly (A:s2 B:s8 C:s8 D:s2 E:s8)
$null = $ws.Range("A44:E44").Copy()
$null = $ws.Range("A46:E46").PasteSpecial(-4122)

# Row 45 = copy row 44 formatting as a base
$null = $ws.Range("A44:E44").Copy()
$null = $ws.Range("A45:E45").PasteSpecial(-4122)

# D45 uses B44's style (s=8) rather than D44's style (s=2)
$null = $ws.Range("B44").Copy()
$null = $ws.Range("D45").PasteSpecial(-4122)

# C45 needs wrap text (new style s=12)
$ws.Range("C45").WrapText = $true

# Values
$ws.Range("A45").Value = "Profile44"
$ws.Range("B45").Value = "OPQA-2936|OPQA-2939|OPQA-2938"
$ws.Range("C45").Value = "Verify that Profile Picture  modal window ‘update’ button should be disabled by default|Verify that Profile picture modal window should contain default buttons and messages|Verify that Profile Picture  modal window should be disappear  while click on ‘Close(X)’ button
"
$ws.Range("D45").Value = "Y"

$ws.Range("A46").Value = "Profile45"
$ws.Range("C46").Value = "Verify that Profile Picture  modal window should be disappear  while click on ‘Cancel’ button"
$ws.Range("B46").Value = "OPQA-2937 "
$ws.Range("D46").Value = "Y"

$ws.Rows.Item(45).RowHeight = 60

# Scroll/selection to match the final view state
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("C55").Select()

Write-Host "done"
